$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-30 down to 18-31
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new market-price entry
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(17, 3).Value = "Los Lagos"
$ws.Cells.Item(17, 4).Value = 44907
$ws.Cells.Item(17, 5).Value = 10
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100101
$ws.Cells.Item(17, 8).Value = "Berries"
$ws.Cells.Item(17, 9).Value = 100101001
$ws.Cells.Item(17, 10).Value = "Arándano (blue)"
$ws.Cells.Item(17, 11).Value = "Sin especificar"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 400
$ws.Cells.Item(17, 14).Value = 4500
$ws.Cells.Item(17, 15).Value = 5000
$ws.Cells.Item(17, 16).Value = 4750
$ws.Cells.Item(17, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(17, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(17, 19).Value = 2375
$ws.Cells.Item(17, 20).Value = 2
